# Updates cryptos list (prices / 1h volume %) per the Fri May 3 09:23:08 UTC 2024
# GitHub Actions scraper run. Price cells that look like plain decimal numbers
# ("566.38", "0.450", "1.00", etc.) are entered with a leading apostrophe so
# Excel stores them as literal text (matching the sheet's existing text-price
# column) instead of auto-converting them to numbers and dropping formatting
# such as trailing zeros. Prices that already contain a thousands-separator
# dot (e.g. "59.385.38") are safe to assign directly since Excel cannot
# parse them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.385.38'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').Value = '2.980.97'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'566.38"
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('D6').Value = "'138.17"
$ws.Range('E6').Value = '  +4.12%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D9').Value = '2.973.01'
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('E10').Value = '  +3.61%  '
$ws.Range('D11').Value = "'5.37"
$ws.Range('E11').Value = '  +11.10%  '
$ws.Range('D12').Value = "'0.450"
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').Value = "'33.64"
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '3.471.55'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = "'7.02"
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '2.978.07'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').Value = '59.400.96'
$ws.Range('E19').Value = '  +2.74%  '
$ws.Range('D20').Value = "'436.18"
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('D21').Value = "'13.57"
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('D22').Value = "'0.720"
$ws.Range('E22').Value = '  +3.15%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = "'79.94"
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = "'2.22"
$ws.Range('E27').Value = '  +10.41%  '
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').Value = "'1.00"
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').Value = "'7.73"
$ws.Range('E30').Value = '  +3.28%  '
$ws.Range('D31').Value = "'25.70"
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').Value = "'6.20"
$ws.Range('E32').Value = '  +4.10%  '
$ws.Range('E33').Value = '  +8.63%  '
$ws.Range('D34').Value = '0.0₃0769'
$ws.Range('E34').Value = '  +10.11%  '
$ws.Range('D35').Value = "'5.89"
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').Value = "'0.982"
$ws.Range('E36').Value = '  +3.92%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').Value = "'48.58"
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = "'8.67"
$ws.Range('E39').Value = '  -3.24%  '
$ws.Range('D40').Value = "'2.77"
$ws.Range('E40').Value = '  +3.50%  '
$ws.Range('D41').Value = "'400.69"
$ws.Range('E41').Value = '  +5.41%  '
$ws.Range('E42').Value = '  +1.11%  '
$ws.Range('D43').Value = '2.732.26'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('E45').Value = '  +5.75%  '
$ws.Range('D47').Value = "'35.03"
$ws.Range('E47').Value = '  +20.26%  '
$ws.Range('D48').Value = "'121.98"
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('D51').Value = "'23.27"
$ws.Range('E51').Value = '  +1.63%  '
